{"js": "// Replace the table-header cell text \"Gerente\" with \"Gerencia\".\n// (The source document contains exactly one run with the text \"Gerente\".)\nconst body = context.document.body;\nconst results = body.search(\"Gerente\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Gerencia\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the table-header cell text \"Gerente\" with \"Gerencia\".\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Execute(\"Gerente\", $true, $true, $false, $false, $false, $true, 1, $false, \"Gerencia\", 2)\n"}
